# regen sval data to filter save games
# Overwrite the numeric stat columns (B:G) for rows 2-20 with the
# regenerated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{}
$data[2]  = @(0.2881169905109251, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 2.598097515653722)
$data[3]  = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 3.755628166162433)
$data[4]  = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 6.15379541431027)
$data[5]  = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 5.582307763322248)
$data[6]  = @(0.6545652718822623, 0.04103571897497393, 0.1496068669990043, 0.5333859586016987, 0, 1.378593816457939)
$data[7]  = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 5.582307763322248)
$data[8]  = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 1, 6.15379541431027)
$data[9]  = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 5.582307763322248)
$data[10] = @(3.272327238179451, 1.626987699542094, 18.71679738969934, 13.86384647080068, 0, 37.47995879822157)
$data[11] = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 1, 6.15379541431027)
$data[12] = @(3.272327238179451, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 0, 4.831698963401555)
$data[13] = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 6.15379541431027)
$data[14] = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 1, 8.656069925401464)
$data[15] = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 3.755628166162433)
$data[16] = @(0.04172184405617529, 0.04103571897497393, 0.7210945179870265, 0.5333859586016987, 1, 1.337238039619874)
$data[17] = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 2797.565817734744, 0, 2803.186227190452)
$data[18] = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 0, 8.656069925401464)
$data[19] = @(0.04172184405617529, 0.0000006194867796516235, 0.7210945179870265, 0.5333859586016987, 0, 1.29620294013168)
$data[20] = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 3.755628166162433)

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
